$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.Style = "Normal"
}

Set-TextValue "D2" "27.938.92"
Set-TextValue "E2" "  +0.60%  "
Set-TextValue "D3" "1.811.83"
Set-TextValue "E3" "  +1.56%  "
Set-TextValue "D4" "1.002"
Set-TextValue "E4" "  -0.09%  "
Set-TextValue "D5" "310.36"
Set-TextValue "E5" "  -0.03%  "
Set-TextValue "D7" "0.4988"
Set-TextValue "E7" "  -2.34%  "
Set-TextValue "D8" "0.3902"
Set-TextValue "E8" "  +1.29%  "
Set-TextValue "D9" "0.09780"
Set-TextValue "E9" "  +25.18%  "
Set-TextValue "E10" "  +0.99%  "
Set-TextValue "D11" "40.84"
Set-TextValue "E11" "  +0.19%  "
Set-TextValue "D12" "6.402"
Set-TextValue "E12" "  +3.37%  "
Set-TextValue "E13" "  -0.12%  "
Set-TextValue "D14" "20.39"
Set-TextValue "E14" "  +1.23%  "
Set-TextValue "D15" "1.812.17"
Set-TextValue "E15" "  +1.81%  "
Set-TextValue "D16" "7.259"
Set-TextValue "E16" "  +0.76%  "
Set-TextValue "D17" "0.00001135"
Set-TextValue "E17" "  +5.70%  "
Set-TextValue "D18" "92.28"
Set-TextValue "E18" "  +1.06%  "
Set-TextValue "D19" "0.06638"
Set-TextValue "E19" "  +1.33%  "
Set-TextValue "D20" "1.001"
Set-TextValue "E20" "  -0.08%  "
Set-TextValue "D21" "17.15"
Set-TextValue "E21" "  +0.89%  "
Set-TextValue "E22" "  -0.09%  "
Set-TextValue "D23" "28.004.14"
Set-TextValue "E23" "  +0.67%  "
Set-TextValue "D24" "11.06"
Set-TextValue "E24" "  +0.65%  "
Set-TextValue "D25" "2.245"
Set-TextValue "E25" "  +0.87%  "
Set-TextValue "D26" "158.28"
Set-TextValue "E26" "  -0.97%  "
Set-TextValue "D27" "2.020.29"
Set-TextValue "E27" "  +1.82%  "
Set-TextValue "D28" "20.52"
Set-TextValue "E28" "  +1.71%  "
Set-TextValue "D29" "2.383"
Set-TextValue "E29" "  +0.83%  "
Set-TextValue "E30" "  +2.67%  "
Set-TextValue "D31" "0.1063"
Set-TextValue "E31" "  -1.14%  "
Set-TextValue "D32" "1.030"
Set-TextValue "E32" "  -0.09%  "
Set-TextValue "D33" "5.552"
Set-TextValue "E33" "  +1.37%  "
Set-TextValue "D34" "3.599"
Set-TextValue "E34" "  -0.98%  "
Set-TextValue "D35" "0.06700"
Set-TextValue "E35" "  -4.96%  "
Set-TextValue "D36" "0.02323"
Set-TextValue "E36" "  +1.01%  "
Set-TextValue "D37" "8.847"
Set-TextValue "E37" "  +0.72%  "
Set-TextValue "E38" "  +0.73%  "
Set-TextValue "D39" "4.926"
Set-TextValue "E39" "  -1.31%  "
Set-TextValue "D40" "11.25"
Set-TextValue "E40" "  -1.52%  "
Set-TextValue "D41" "0.6158"
Set-TextValue "E41" "  +1.36%  "
Set-TextValue "D42" "1.176"
Set-TextValue "E42" "  +2.40%  "
Set-TextValue "E43" "  -0.10%  "
Set-TextValue "D44" "13.18"
Set-TextValue "E44" "  +0.79%  "
Set-TextValue "D45" "0.5882"
Set-TextValue "E45" "  -0.15%  "
Set-TextValue "B46" "PancakeSwap"
Set-TextValue "C46" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D46" "3.692"
Set-TextValue "E46" "  -0.23%  "
Set-TextValue "B47" "WEMIXTOKEN"
Set-TextValue "C47" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D47" "1.281"
Set-TextValue "E47" "  -2.68%  "
Set-TextValue "D48" "123.50"
Set-TextValue "E48" "  -1.85%  "
Set-TextValue "E49" "  +1.93%  "
Set-TextValue "E50" "  -1.86%  "
Set-TextValue "D51" "0.06768"
Set-TextValue "E51" "  -1.05%  "
